$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("postBatchValid")
$ws.Activate()

# Update the duplicate step-def values on row 4
$ws.Range("B4").Value = "Salesforce1"
$ws.Range("C4").Value = "Active"

# Move the active selection from C4 to B4
$ws.Range("B4").Select()
